$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: add "commercial" variety in column C (first new shared string)
$ws.Range("C18").Value = "commercial"

# Add new "Price" column header (second new shared string)
$ws.Range("F1").Value = "Price"

# Fill price of 2.99 for all data rows 2-24
for ($r = 2; $r -le 24; $r++) {
    $ws.Cells.Item($r, 6).Value = 2.99
}

# Row 24: fill in missing Veg type (A24) and variety "Scotian" (C24, third new shared string)
$ws.Range("A24").Value = "Veg"
$ws.Range("C24").Value = "Scotian"

# Update selection to F24 to match final state
$ws.Range("F24").Select()
